# Change the table style used by the "Sources of finance" table on slide 6
# from its original style GUID to the new one.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$tbl.ApplyStyle("{2DB7A4E9-1E38-4EFF-B8D3-8BB6C8B1F3FA}", $true)
